$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "magasin" sheet: fill in the magasin data table (B2:F6), a stray note
#    cell (D8), fix the lienInstrumentMagasin formula casing (D11:D61),
#    and refresh the selection / column widths.
# ---------------------------------------------------------------------------
$wsMagasin = $wb.Worksheets.Item("magasin")

$wsMagasin.Range("B2").Value = "26 Av. des Bénédictins"
$wsMagasin.Range("C2").Value = 87000
$wsMagasin.Range("D2").Value = "Le magasin de musique de Limoges. Une large gamme de guitares acoustiques, électriques, électro-ac mais aussi de basses, pianos, claviers et batteries."
$wsMagasin.Range("E2").Value = "Music Passion 87"
$wsMagasin.Range("F2").Value = "45.832950, 1.267190"

$wsMagasin.Range("B3").Value = "8 Rue des Combes BP 30190"
$wsMagasin.Range("C3").Value = 87005
$wsMagasin.Range("D3").Value = "Située dans le centre-ville de Limoges, la Fnac vous accueille dans son magasin d'une superficie de plus de 2100m2. Vous y retrouverez tous les univers de la Fnac : livres, disques, DVD, micro‐informatique, objets connectés, gaming, téléphonie, son, photo, TV, papeterie,… Pour vos places de spectacles l'espace billetterie de Fnac Spectacles propose une large sélection d'événements."
$wsMagasin.Range("E3").Value = "FNAC Limoges"
$wsMagasin.Range("F3").Value = "45.832250, 1.257570"

$wsMagasin.Range("B4").Value = "30 Rue Amédée Gordini"
$wsMagasin.Range("C4").Value = 87280
$wsMagasin.Range("D4").Value = "Cultura enseigne leader de biens de loisirs culturels et créatifs"
$wsMagasin.Range("E4").Value = "Cultura Limoges"
$wsMagasin.Range("F4").Value = "45.895320, 1.280380"

$wsMagasin.Range("B5").Value = "12 Rue Jules Guesde"
$wsMagasin.Range("C5").Value = 87000
$wsMagasin.Range("D5").Value = "Petit magasin proposant un large eventail de produits"
$wsMagasin.Range("E5").Value = "Music Mania"
$wsMagasin.Range("F5").Value = "45.829520, 1.261290"

$wsMagasin.Range("B6").Value = "5 Rue de la Glâne"
$wsMagasin.Range("C6").Value = 87000
$wsMagasin.Range("D6").Value = "rien"
$wsMagasin.Range("E6").Value = "Limouzik"
$wsMagasin.Range("F6").Value = "45.822910, 1.218190"

$wsMagasin.Range("D8").Value = "5 Rue de la Glâne, 87000 Limoges"

$wsMagasin.Range("D11").Formula = '="INSERT INTO lienInstrumentMagasin (idInstrument, idMagasin, stock) VALUES("&A11&","&B11&","&C11&");"'
$wsMagasin.Range("D12:D61").Formula = '="INSERT INTO lienInstrumentMagasin (idInstrument, idMagasin, stock) VALUES("&A12&","&B12&","&C12&");"'

$wsMagasin.Columns.Item(2).AutoFit()
$wsMagasin.Columns.Item(6).AutoFit()

$wsMagasin.Range("G2:G6").Select()
$wsMagasin.Range("G2").Select()

# ---------------------------------------------------------------------------
# 2) "table lien couleur" sheet: header label "nomimage" -> "nomImage".
# ---------------------------------------------------------------------------
$wsCouleur = $wb.Worksheets.Item("table lien couleur")
$wsCouleur.Range("E21").Value = "nomImage"
$wsCouleur.Range("C28").Select()

# ---------------------------------------------------------------------------
# 3) "Table instru" sheet becomes the active / selected tab.
# ---------------------------------------------------------------------------
$wsInstru = $wb.Worksheets.Item("Table instru")
$wsInstru.Activate()
$wsInstru.Range("C26").Select()
$excel.ActiveWindow.Zoom = 85
